# Atualização automática via cronjob
# Atualiza a base de vendas atípicas: corrige quantidades de estoque
# recalculadas, renumera os dias reprocessados e inclui a nova ocorrência
# detectada na varredura mais recente.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrige estoque_atualizado (coluna G) para registros recalculados ---
$ws.Range("G2").Value = 9
$ws.Range("G3").Value = -40
$ws.Range("G8").Value = 94
$ws.Range("G9").Value = 803

# --- Renumeração do contador (coluna A) apos inclusao da nova linha ---
$ws.Range("A6").Value = 7
$ws.Range("A7").Value = 8

# --- Nova ocorrencia atipica detectada (linha 10) ---
$ws.Range("B10").Value = "'2025-04-22"
$ws.Range("B10").ClearFormats()

$ws.Range("C10").Value = 24

$ws.Range("D10").Value = "BRAGA MOTOS LTDA"

$ws.Range("E10").Value = "'000015"
$ws.Range("E10").ClearFormats()

$ws.Range("F10").Value = "PANO MULTIUSO ROLO 28X300 M AZUL TALGE"

$ws.Range("G10").Value = 7

$ws.Range("H10").Value = $false

# Replica o estilo (borda/negrito) usado nas demais celulas da coluna A
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A10").Value = 6
